# Auto-generated Excel COM-interop script
# Refreshes cached market-board price / leve-profit figures across all
# 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the
# latest scrape from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 2974.353
$ws.Range("J17").Value2 = 3290.2856
$ws.Range("L17").Value2 = 9870.856800000001
$ws.Range("N17").Value2 = -10206.8568

$ws.Range("H62").Value2 = 4509.2104
$ws.Range("I62").Value2 = 3262.5
$ws.Range("J62").Value2 = 8000
$ws.Range("K62").Value2 = 3262.5
$ws.Range("L62").Value2 = 8000
$ws.Range("M62").Value2 = -2638.5
$ws.Range("N62").Value2 = -9248

$ws.Range("H65").Value2 = 4509.2104
$ws.Range("I65").Value2 = 3262.5
$ws.Range("J65").Value2 = 8000
$ws.Range("K65").Value2 = 16312.5
$ws.Range("L65").Value2 = 40000
$ws.Range("M65").Value2 = -13192.5
$ws.Range("N65").Value2 = -46240

$ws.Range("I88").Value2 = 1086
$ws.Range("J88").Value2 = 5502
$ws.Range("K88").Value2 = 1086
$ws.Range("L88").Value2 = 5502
$ws.Range("M88").Value2 = -680
$ws.Range("N88").Value2 = -6314

$ws.Range("I91").Value2 = 1086
$ws.Range("J91").Value2 = 5502
$ws.Range("K91").Value2 = 1086
$ws.Range("L91").Value2 = 5502
$ws.Range("M91").Value2 = 318
$ws.Range("N91").Value2 = -8310

$ws.Range("H125").Value2 = 745.5
$ws.Range("J125").Value2 = 444.5
$ws.Range("L125").Value2 = 4000.5
$ws.Range("N125").Value2 = -8920.5

$ws.Range("H129").Value2 = 4160.75
$ws.Range("J129").Value2 = 5800
$ws.Range("L129").Value2 = 17400
$ws.Range("N129").Value2 = -27400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("K4").Value2 = 0
$ws.Range("M4").ClearContents()

$ws.Range("H5").Value2 = 680
$ws.Range("I5").Value2 = 680
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 680
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = -568
$ws.Range("N5").ClearContents()

$ws.Range("H32").Value2 = 18272.283
$ws.Range("I32").Value2 = 8891.68
$ws.Range("J32").Value2 = 26647.822
$ws.Range("K32").Value2 = 8891.68
$ws.Range("L32").Value2 = 26647.822
$ws.Range("M32").Value2 = -8604.68
$ws.Range("N32").Value2 = -27221.822

$ws.Range("H88").Value2 = 1900
$ws.Range("I88").Value2 = 1656
$ws.Range("J88").Value2 = 2192.8
$ws.Range("K88").Value2 = 1656
$ws.Range("L88").Value2 = 2192.8
$ws.Range("M88").Value2 = -1250
$ws.Range("N88").Value2 = -3004.8

$ws.Range("H91").Value2 = 1900
$ws.Range("I91").Value2 = 1656
$ws.Range("J91").Value2 = 2192.8
$ws.Range("K91").Value2 = 1656
$ws.Range("L91").Value2 = 2192.8
$ws.Range("M91").Value2 = -252
$ws.Range("N91").Value2 = -5000.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 680
$ws.Range("I4").Value2 = 680
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 680
$ws.Range("L4").Value2 = 0
$ws.Range("M4").Value2 = -565
$ws.Range("N4").ClearContents()

$ws.Range("H20").Value2 = 1650
$ws.Range("I20").Value2 = 1533.3334
$ws.Range("J20").Value2 = 2000
$ws.Range("K20").Value2 = 1533.3334
$ws.Range("L20").Value2 = 2000
$ws.Range("M20").Value2 = -1286.3334
$ws.Range("N20").Value2 = -2494

$ws.Range("H36").Value2 = 965.3333
$ws.Range("I36").Value2 = 949.5
$ws.Range("J36").Value2 = 997
$ws.Range("K36").Value2 = 949.5
$ws.Range("L36").Value2 = 997
$ws.Range("M36").Value2 = -415.5
$ws.Range("N36").Value2 = -2065

$ws.Range("H75").Value2 = 50398.8
$ws.Range("I75").Value2 = 12000
$ws.Range("K75").Value2 = 12000
$ws.Range("M75").Value2 = -11064

$ws.Range("H78").Value2 = 50398.8
$ws.Range("I78").Value2 = 12000
$ws.Range("K78").Value2 = 36000
$ws.Range("M78").Value2 = -31320

$ws.Range("H99").Value2 = 1701.8334
$ws.Range("I99").Value2 = 1543.2
$ws.Range("K99").Value2 = 1543.2
$ws.Range("M99").Value2 = -45.20000000000005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 5315.2
$ws.Range("I31").Value2 = 2811.5
$ws.Range("J31").Value2 = 6984.3335
$ws.Range("K31").Value2 = 2811.5
$ws.Range("L31").Value2 = 6984.3335
$ws.Range("M31").Value2 = -2516.5
$ws.Range("N31").Value2 = -7574.3335

$ws.Range("H34").Value2 = 5315.2
$ws.Range("I34").Value2 = 2811.5
$ws.Range("J34").Value2 = 6984.3335
$ws.Range("K34").Value2 = 2811.5
$ws.Range("L34").Value2 = 6984.3335
$ws.Range("M34").Value2 = -2609.5
$ws.Range("N34").Value2 = -7388.3335

$ws.Range("H62").Value2 = 32362.857
$ws.Range("I62").Value2 = 3865.4443
$ws.Range("K62").Value2 = 3865.4443
$ws.Range("M62").Value2 = -3241.4443

$ws.Range("H65").Value2 = 32362.857
$ws.Range("I65").Value2 = 3865.4443
$ws.Range("K65").Value2 = 19327.2215
$ws.Range("M65").Value2 = -16207.2215

$ws.Range("H86").Value2 = 10738.429
$ws.Range("I86").Value2 = 7715.625
$ws.Range("K86").Value2 = 7715.625
$ws.Range("M86").Value2 = -6592.625

$ws.Range("H89").Value2 = 10738.429
$ws.Range("I89").Value2 = 7715.625
$ws.Range("K89").Value2 = 38578.125
$ws.Range("M89").Value2 = -32962.125

$ws.Range("H132").Value2 = 2887.0688
$ws.Range("I132").Value2 = 2708.1428
$ws.Range("K132").Value2 = 8124.428400000001
$ws.Range("M132").Value2 = -5594.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value2 = 403
$ws.Range("I92").Value2 = 303.75
$ws.Range("J92").Value2 = 800
$ws.Range("K92").Value2 = 911.25
$ws.Range("L92").Value2 = 2400
$ws.Range("M92").Value2 = 336.75
$ws.Range("N92").Value2 = -4896

$ws.Range("H131").Value2 = 3112.6924
$ws.Range("J131").Value2 = 3521
$ws.Range("L131").Value2 = 10563
$ws.Range("N131").Value2 = -20643

$ws.Range("H139").Value2 = 3144.8235
$ws.Range("I139").Value2 = 2622
$ws.Range("K139").Value2 = 7866
$ws.Range("M139").Value2 = -2726

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 8109.778
$ws.Range("J70").Value2 = 8247.5
$ws.Range("L70").Value2 = 8247.5
$ws.Range("N70").Value2 = -8787.5

$ws.Range("H73").Value2 = 8109.778
$ws.Range("J73").Value2 = 8247.5
$ws.Range("L73").Value2 = 8247.5
$ws.Range("N73").Value2 = -10119.5

$ws.Range("H102").Value2 = 1382.5758
$ws.Range("I102").Value2 = 757
$ws.Range("K102").Value2 = 757
$ws.Range("M102").Value2 = 865

$ws.Range("H122").Value2 = 689917.7
$ws.Range("I122").Value2 = 85224.5
$ws.Range("K122").Value2 = 255673.5
$ws.Range("M122").Value2 = -253223.5

$ws.Range("H126").Value2 = 4312.5
$ws.Range("I126").Value2 = 2750
$ws.Range("J126").Value2 = 4833.3335
$ws.Range("K126").Value2 = 8250
$ws.Range("L126").Value2 = 14500.0005
$ws.Range("M126").Value2 = -5780
$ws.Range("N126").Value2 = -19440.0005

$ws.Range("H136").Value2 = 64999.75
$ws.Range("J136").Value2 = 64999.75
$ws.Range("L136").Value2 = 194999.25
$ws.Range("N136").Value2 = -200099.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 940.2857
$ws.Range("J22").Value2 = 816.4
$ws.Range("L22").Value2 = 816.4
$ws.Range("N22").Value2 = -1406.4

$ws.Range("H27").Value2 = 940.2857
$ws.Range("J27").Value2 = 816.4
$ws.Range("L27").Value2 = 816.4
$ws.Range("N27").Value2 = -1030.4

$ws.Range("H40").Value2 = 3126.5
$ws.Range("I40").Value2 = 3002
$ws.Range("J40").Value2 = 3500
$ws.Range("K40").Value2 = 3002
$ws.Range("L40").Value2 = 3500
$ws.Range("M40").Value2 = -2866
$ws.Range("N40").Value2 = -3772

$ws.Range("H122").Value2 = 3663
$ws.Range("I122").Value2 = 3649.5
$ws.Range("J122").Value2 = 3690
$ws.Range("K122").Value2 = 10948.5
$ws.Range("L122").Value2 = 11070
$ws.Range("M122").Value2 = -8498.5
$ws.Range("N122").Value2 = -15970

$ws.Range("H132").Value2 = 3097.647
$ws.Range("I132").Value2 = 2032.8
$ws.Range("J132").Value2 = 6055.5557
$ws.Range("K132").Value2 = 6098.4
$ws.Range("L132").Value2 = 18166.6671
$ws.Range("M132").Value2 = -3568.4
$ws.Range("N132").Value2 = -23226.6671

$ws.Range("H136").Value2 = 2997.6667
$ws.Range("I136").Value2 = 2997.6667
$ws.Range("K136").Value2 = 8993.000100000001
$ws.Range("M136").Value2 = -6443.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 7905.9165
$ws.Range("I62").Value2 = 6985.5
$ws.Range("K62").Value2 = 6985.5
$ws.Range("M62").Value2 = -6361.5

$ws.Range("H65").Value2 = 7905.9165
$ws.Range("I65").Value2 = 6985.5
$ws.Range("K65").Value2 = 34927.5
$ws.Range("M65").Value2 = -31807.5

$ws.Range("H81").Value2 = 1450.08
$ws.Range("J81").Value2 = 0
$ws.Range("L81").Value2 = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value2 = 1450.08
$ws.Range("J84").Value2 = 0
$ws.Range("L84").Value2 = 0
$ws.Range("N84").ClearContents()

$ws.Range("H113").Value2 = 3429.7693
$ws.Range("I113").Value2 = 2475.25
$ws.Range("J113").Value2 = 3854
$ws.Range("K113").Value2 = 7425.75
$ws.Range("L113").Value2 = 11562
$ws.Range("M113").Value2 = -5255.75
$ws.Range("N113").Value2 = -15902

$ws.Range("H141").Value2 = 50000
$ws.Range("J141").Value2 = 50000
$ws.Range("L141").Value2 = 50000
$ws.Range("N141").Value2 = -60360
